$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-05 Monday", "2025-05-06 Tuesday"),
    @("83×94=7802", "34×32=1088"),
    @("20×80=1600", "79×25=1975"),
    @("65×80=5200", "37×84=3108"),
    @("63×41=2583", "54×61=3294"),
    @("34×98=3332", "65×95=6175"),
    @("84×72=6048", "97×57=5529"),
    @("12×52=624", "63×21=1323"),
    @("45×84=3780", "49×27=1323"),
    @("98×90=8820", "85×15=1275"),
    @("29×48=1392", "98×87=8526"),
    @("36×60=2160", "87×14=1218"),
    @("97×34=3298", "55×11=605"),
    @("25×82=2050", "72×87=6264"),
    @("65×58=3770", "49×25=1225"),
    @("82×63=5166", "46×52=2392"),
    @("78×84=6552", "92×13=1196"),
    @("52×92=4784", "46×65=2990"),
    @("39×95=3705", "15×63=945"),
    @("69×60=4140", "58×60=3480"),
    @("22×99=2178", "19×85=1615"),
    @("11×25=275", "54×41=2214"),
    @("63×72=4536", "22×79=1738"),
    @("21×70=1470", "12×45=540"),
    @("49×67=3283", "55×82=4510"),
    @("92×55=5060", "20×28=560")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
